$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "26.480.78"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "1.842.46"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "261.51"
$ws.Range("E5").Value = "  -5.73%  "
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.5190"
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("D8").Value = "0.3281"
$ws.Range("E8").Value = "  -4.08%  "
$ws.Range("D9").Value = "0.06796"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "18.72"
$ws.Range("E10").Value = "  -6.53%  "
$ws.Range("D11").Value = "0.7712"
$ws.Range("E11").Value = "  -3.98%  "
$ws.Range("D12").Value = "0.07699"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "1.807.01"
$ws.Range("E13").Value = "  -4.07%  "
$ws.Range("D14").Value = "88.40"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").Value = "0.9990"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "13.94"
$ws.Range("E17").Value = "  -4.19%  "
$ws.Range("D18").Value = "0.000007989"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "26.472.62"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").Value = "2.068.18"
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").Value = "4.579"
$ws.Range("E22").Value = "  -3.39%  "
$ws.Range("D23").Value = "9.490"
$ws.Range("E23").Value = "  -5.29%  "
$ws.Range("D24").Value = "5.973"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").Value = "144.23"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("D26").Value = "2.219"
$ws.Range("E26").Value = "  -7.92%  "
$ws.Range("D27").Value = "1.647"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "17.01"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("D29").Value = "111.63"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "4.199"
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("D31").Value = "4.144"
$ws.Range("E31").Value = "  -3.65%  "
$ws.Range("D32").Value = "0.08732"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("D33").Value = "0.04812"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").Value = "1.134"
$ws.Range("E34").Value = "  -3.51%  "
$ws.Range("D35").Value = "2.838"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("D36").Value = "0.7080"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "3.079"
$ws.Range("E37").Value = "  -6.38%  "
$ws.Range("D38").Value = "2.239"
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("D39").Value = "0.01760"
$ws.Range("E39").Value = "  -4.38%  "
$ws.Range("D40").Value = "0.4854"
$ws.Range("E40").Value = "  -5.14%  "
$ws.Range("D41").Value = "111.51"
$ws.Range("E41").Value = "  -4.04%  "
$ws.Range("D42").Value = "0.8943"
$ws.Range("E42").Value = "  -6.38%  "
$ws.Range("D43").Value = "6.082"
$ws.Range("D44").Value = "0.9992"
$ws.Range("D45").Value = "7.730"
$ws.Range("E45").Value = "  -4.42%  "
$ws.Range("D46").Value = "0.4156"
$ws.Range("E46").Value = "  -6.63%  "
$ws.Range("D47").Value = "0.05870"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").Value = "8.975"
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("D49").Value = "35.04"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "0.1221"
$ws.Range("E50").Value = "  -8.72%  "
$ws.Range("E51").Value = "  +0.80%  "
